# Update LR-pair table (rows 2-4 modified, rows 5-10 added)
# per "Natmi following Dr Hou advice" commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value2 = "ECs"
$ws.Range("B2").Value2 = "Calr"
$ws.Range("C2").Value2 = "Itga3"
$ws.Range("D2").Value2 = "ECs"
$ws.Range("E2").Value2 = 3
$ws.Range("F2").Value2 = 1
$ws.Range("G2").Value2 = 72.79331566666666
$ws.Range("H2").Value2 = 218.379947
$ws.Range("I2").Value2 = 0.2828741606141505
$ws.Range("J2").Value2 = 0.2828741606141506
$ws.Range("K2").Value2 = 2
$ws.Range("L2").Value2 = 0.6666666666666666
$ws.Range("M2").Value2 = 6.169512999999999
$ws.Range("N2").Value2 = 18.508539
$ws.Range("O2").Value2 = 0.5207942167525852
$ws.Range("P2").Value2 = 0.5207942167525853
$ws.Range("Q2").Value2 = 449.0993073186036
$ws.Range("R2").Value2 = 4041.893765867432
$ws.Range("S2").Value2 = 0.1473192269165915
$ws.Range("T2").Value2 = 0.1473192269165916

# Row 3
$ws.Range("A3").Value2 = "ECs"
$ws.Range("B3").Value2 = "Calr"
$ws.Range("C3").Value2 = "Itga3"
$ws.Range("D3").Value2 = "FAPs"
$ws.Range("E3").Value2 = 3
$ws.Range("F3").Value2 = 1
$ws.Range("G3").Value2 = 72.79331566666666
$ws.Range("H3").Value2 = 218.379947
$ws.Range("I3").Value2 = 0.2828741606141505
$ws.Range("J3").Value2 = 0.2828741606141506
$ws.Range("K3").Value2 = 2
$ws.Range("L3").Value2 = 0.6666666666666666
$ws.Range("M3").Value2 = 0.06813733333333333
$ws.Range("N3").Value2 = 0.204412
$ws.Range("O3").Value2 = 0.005751755307905689
$ws.Range("P3").Value2 = 0.00575175530790569
$ws.Range("Q3").Value2 = 4.959942414018221
$ws.Range("R3").Value2 = 44.63948172616399
$ws.Range("S3").Value2 = 0.001627022954781807
$ws.Range("T3").Value2 = 0.001627022954781807

# Row 4
$ws.Range("A4").Value2 = "ECs"
$ws.Range("B4").Value2 = "Calr"
$ws.Range("C4").Value2 = "Itga3"
$ws.Range("D4").Value2 = "sCs"
$ws.Range("E4").Value2 = 3
$ws.Range("F4").Value2 = 1
$ws.Range("G4").Value2 = 72.79331566666666
$ws.Range("H4").Value2 = 218.379947
$ws.Range("I4").Value2 = 0.2828741606141505
$ws.Range("J4").Value2 = 0.2828741606141506
$ws.Range("K4").Value2 = 3
$ws.Range("L4").Value2 = 1
$ws.Range("M4").Value2 = 5.608704333333333
$ws.Range("N4").Value2 = 16.826113
$ws.Range("O4").Value2 = 0.473454027939509
$ws.Range("P4").Value2 = 0.4734540279395091
$ws.Range("Q4").Value2 = 408.2761850173345
$ws.Range("R4").Value2 = 3674.48566515601
$ws.Range("S4").Value2 = 0.1339279107427772
$ws.Range("T4").Value2 = 0.1339279107427772

# Row 5
$ws.Range("A5").Value2 = "FAPs"
$ws.Range("B5").Value2 = "Calr"
$ws.Range("C5").Value2 = "Itga3"
$ws.Range("D5").Value2 = "ECs"
$ws.Range("E5").Value2 = 3
$ws.Range("F5").Value2 = 1
$ws.Range("G5").Value2 = 117.1700846666667
$ws.Range("H5").Value2 = 351.510254
$ws.Range("I5").Value2 = 0.4553218801152877
$ws.Range("J5").Value2 = 0.4553218801152878
$ws.Range("K5").Value2 = 2
$ws.Range("L5").Value2 = 0.6666666666666666
$ws.Range("M5").Value2 = 6.169512999999999
$ws.Range("N5").Value2 = 18.508539
$ws.Range("O5").Value2 = 0.5207942167525852
$ws.Range("P5").Value2 = 0.5207942167525853
$ws.Range("Q5").Value2 = 722.8823605621006
$ws.Range("R5").Value2 = 6505.941245058906
$ws.Range("S5").Value2 = 0.2371290019249557
$ws.Range("T5").Value2 = 0.2371290019249558

# Row 6
$ws.Range("A6").Value2 = "FAPs"
$ws.Range("B6").Value2 = "Calr"
$ws.Range("C6").Value2 = "Itga3"
$ws.Range("D6").Value2 = "FAPs"
$ws.Range("E6").Value2 = 3
$ws.Range("F6").Value2 = 1
$ws.Range("G6").Value2 = 117.1700846666667
$ws.Range("H6").Value2 = 351.510254
$ws.Range("I6").Value2 = 0.4553218801152877
$ws.Range("J6").Value2 = 0.4553218801152878
$ws.Range("K6").Value2 = 2
$ws.Range("L6").Value2 = 0.6666666666666666
$ws.Range("M6").Value2 = 0.06813733333333333
$ws.Range("N6").Value2 = 0.204412
$ws.Range("O6").Value2 = 0.005751755307905689
$ws.Range("P6").Value2 = 0.00575175530790569
$ws.Range("Q6").Value2 = 7.983657115627556
$ws.Range("R6").Value2 = 71.852914040648
$ws.Range("S6").Value2 = 0.002618900040758704
$ws.Range("T6").Value2 = 0.002618900040758705

# Row 7
$ws.Range("A7").Value2 = "FAPs"
$ws.Range("B7").Value2 = "Calr"
$ws.Range("C7").Value2 = "Itga3"
$ws.Range("D7").Value2 = "sCs"
$ws.Range("E7").Value2 = 3
$ws.Range("F7").Value2 = 1
$ws.Range("G7").Value2 = 117.1700846666667
$ws.Range("H7").Value2 = 351.510254
$ws.Range("I7").Value2 = 0.4553218801152877
$ws.Range("J7").Value2 = 0.4553218801152878
$ws.Range("K7").Value2 = 3
$ws.Range("L7").Value2 = 1
$ws.Range("M7").Value2 = 5.608704333333333
$ws.Range("N7").Value2 = 16.826113
$ws.Range("O7").Value2 = 0.473454027939509
$ws.Range("P7").Value2 = 0.4734540279395091
$ws.Range("Q7").Value2 = 657.1723616069669
$ws.Range("R7").Value2 = 5914.551254462703
$ws.Range("S7").Value2 = 0.2155739781495732
$ws.Range("T7").Value2 = 0.2155739781495733

# Row 8
$ws.Range("A8").Value2 = "sCs"
$ws.Range("B8").Value2 = "Calr"
$ws.Range("C8").Value2 = "Itga3"
$ws.Range("D8").Value2 = "ECs"
$ws.Range("E8").Value2 = 3
$ws.Range("F8").Value2 = 1
$ws.Range("G8").Value2 = 67.37122333333333
$ws.Range("H8").Value2 = 202.11367
$ws.Range("I8").Value2 = 0.2618039592705617
$ws.Range("J8").Value2 = 0.2618039592705618
$ws.Range("K8").Value2 = 2
$ws.Range("L8").Value2 = 0.6666666666666666
$ws.Range("M8").Value2 = 6.169512999999999
$ws.Range("N8").Value2 = 18.508539
$ws.Range("O8").Value2 = 0.5207942167525852
$ws.Range("P8").Value2 = 0.5207942167525853
$ws.Range("Q8").Value2 = 415.6476381809033
$ws.Range("R8").Value2 = 3740.82874362813
$ws.Range("S8").Value2 = 0.1363459879110379
$ws.Range("T8").Value2 = 0.136345987911038

# Row 9
$ws.Range("A9").Value2 = "sCs"
$ws.Range("B9").Value2 = "Calr"
$ws.Range("C9").Value2 = "Itga3"
$ws.Range("D9").Value2 = "FAPs"
$ws.Range("E9").Value2 = 3
$ws.Range("F9").Value2 = 1
$ws.Range("G9").Value2 = 67.37122333333333
$ws.Range("H9").Value2 = 202.11367
$ws.Range("I9").Value2 = 0.2618039592705617
$ws.Range("J9").Value2 = 0.2618039592705618
$ws.Range("K9").Value2 = 2
$ws.Range("L9").Value2 = 0.6666666666666666
$ws.Range("M9").Value2 = 0.06813733333333333
$ws.Range("N9").Value2 = 0.204412
$ws.Range("O9").Value2 = 0.005751755307905689
$ws.Range("P9").Value2 = 0.00575175530790569
$ws.Range("Q9").Value2 = 4.590495501337777
$ws.Range("R9").Value2 = 41.31445951204
$ws.Range("S9").Value2 = 0.001505832312365178
$ws.Range("T9").Value2 = 0.001505832312365179

# Row 10
$ws.Range("A10").Value2 = "sCs"
$ws.Range("B10").Value2 = "Calr"
$ws.Range("C10").Value2 = "Itga3"
$ws.Range("D10").Value2 = "sCs"
$ws.Range("E10").Value2 = 3
$ws.Range("F10").Value2 = 1
$ws.Range("G10").Value2 = 67.37122333333333
$ws.Range("H10").Value2 = 202.11367
$ws.Range("I10").Value2 = 0.2618039592705617
$ws.Range("J10").Value2 = 0.2618039592705618
$ws.Range("K10").Value2 = 3
$ws.Range("L10").Value2 = 1
$ws.Range("M10").Value2 = 5.608704333333333
$ws.Range("N10").Value2 = 16.826113
$ws.Range("O10").Value2 = 0.473454027939509
$ws.Range("P10").Value2 = 0.4734540279395091
$ws.Range("Q10").Value2 = 377.8652722516344
$ws.Range("R10").Value2 = 3400.78745026471
$ws.Range("S10").Value2 = 0.1239521390471586
$ws.Range("T10").Value2 = 0.1239521390471587
